# Sprint 5 backlog update: mark all "In Progress" items as "Complete"
# and move the active selection to A11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Progress column (A2:A7) for all backlog rows from
# "In Progress" to "Complete".
$ws.Range("A2:A7").Value = "Complete"

# Reflect the updated selection/active cell saved in the sheet view.
$ws.Range("A11").Select()
